$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SpecialPrices")

# Fill in the new rows 5-7 (phone/sku/special_price columns A-C)
$ws.Range("A5").Value = 70333029
$ws.Range("B5").Value = "sku-3"
$ws.Range("C5").Value = 2

$ws.Range("A6").Value = 70333029
$ws.Range("B6").Value = "sku-4"
$ws.Range("C6").Value = 2.5

$ws.Range("A7").Value = 70333029
$ws.Range("B7").Value = "sku-5"
$ws.Range("C7").Value = 1.5

# Update the selected cell to match the author's final cursor position
$ws.Range("C6").Select()

$wb.Save()
